$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 124, shifting existing rows 124:241 down to 125:242.
$ws.Rows("124:124").Insert()

# Populate the newly inserted row 124 with its data (weekly entry added between
# the old row 123 and what is now row 125).
$ws.Range("A124").Value = 5
$ws.Range("B124").Value = "Macroferia Regional de Talca"
$ws.Range("C124").Value = "Maule"
$ws.Range("D124").Value = 44705
$ws.Range("E124").Value = 7
$ws.Range("F124").Value = 100112008
$ws.Range("G124").Value = "Coliflor"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 4000
$ws.Range("K124").Value = 900
$ws.Range("L124").Value = 900
$ws.Range("M124").Value = 900
$ws.Range("N124").Value = "$/unidad"
$ws.Range("O124").Value = "Región del Maule"
$ws.Range("P124").Value = 900
$ws.Range("Q124").Value = 1
$ws.Range("R124").Value = "Hortaliza"
